# Fix ART-DECOR server URL in the Codebooks sheet: the DECOR services
# endpoint moved from https://decor.nictiz.nl/services/ to
# https://decor.nictiz.nl/decor/services/. Update every ARTDECOR codebook
# row (rows 2-12, column D "server") to the corrected URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebooks")

# Column D ("server") for every ARTDECOR entry (rows 2-12) held the old URL.
$ws.Range("D2:D12").Value = "https://decor.nictiz.nl/decor/services/"

# Restore the sheet's active selection to D2 (matches the author's saved
# cursor position after making the edit).
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
